# Generate Report for Handoff
# - Flip Status from "In Translation" to "Ready for handoff" on every sheet
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Columns holding the Status text grow wider to fit the new, longer label

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("G2").Value = "2016-08-23 15:14:22"
$ovw.Columns("E:E").ColumnWidth = 16.333333
$ovw.Columns("F:F").ColumnWidth = 16.333333

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-23 15:14:11"
$zhcn.Columns("C:C").ColumnWidth = 16.333333

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-23 15:14:22"
$dede.Columns("C:C").ColumnWidth = 16.333333
